$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (dates are Excel serial numbers, same style as prior rows)
$data = @(
    @(44319, 3, 30, 297.914597815293),
    @(44320, 0, 29, 287.9841112214499),
    @(44321, 0, 29, 287.9841112214499)
)

$startRow = 245
$lastExistingRow = $startRow - 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]

    # Copy the date-cell formatting from the previous row so the new cell
    # gets the same style (centered, bordered, date-formatted) used for
    # all the other rows in column A.
    $ws.Cells.Item($lastExistingRow, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false
